$wb = $excel.ActiveWorkbook

# 1. Update shared text "Ready for handoff" -> "Handoff transform failed"
#    (every cell that currently shows "Ready for handoff" needs updating so the
#     shared string itself is replaced rather than leaving an orphaned string)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value() = "Handoff transform failed"
$wsOverview.Range("C2").Value() = "Handoff transform failed"

# 2. Update zh-cn and de-de sheets
$sheetNames = @("zh-cn", "de-de")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("B2").Value() = "Handoff transform failed"

    # remove hyperlink + cell at C2 (Latest Handoff File)
    $target = $ws.Range("C2").Address()
    foreach ($hl in @($ws.Hyperlinks)) {
        $addr = $hl.Range.Address()
        if ($addr -eq $target) {
            $hl.Delete()
        }
    }
    $ws.Range("C2").Clear()

    # update D2 (Latest Handoff Datetime) to the zero date
    $ws.Range("D2").Value() = "0001-01-01 00:00:00"

    # update H2 (Handoff Reason) Include -> Ignored
    $ws.Range("H2").Value() = "Ignored"
}

Write-Output "Done"
